# The commit adds one new weekly price observation for Jengibre
# (Mercado Mayorista Lo Valledor de Santiago) as a new row inserted at
# row 3 — every existing record from the old row 3 onward shifts down
# by one row (dimension grows from A1:R108 to A1:R109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 3..108 down to 4..109, leaving a blank row 3
# (Excel copies formatting from the row above, same as a real
# right-click "Insert" in the UI).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new observation.
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(3, 3).Value = 'Metropolitana'
$ws.Cells.Item(3, 4).Value = 44860
$ws.Cells.Item(3, 5).Value = 13
$ws.Cells.Item(3, 6).Value = 100114007
$ws.Cells.Item(3, 7).Value = 'Jengibre'
$ws.Cells.Item(3, 8).Value = 'Sin especificar'
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 580
$ws.Cells.Item(3, 11).Value = 11000
$ws.Cells.Item(3, 12).Value = 12000
$ws.Cells.Item(3, 13).Value = 11448
$ws.Cells.Item(3, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(3, 15).Value = 'Perú'
$ws.Cells.Item(3, 16).Value = 881
$ws.Cells.Item(3, 17).Value = 13
$ws.Cells.Item(3, 18).Value = 'Hortaliza'
